# "fix nullable in plan visit"
# Rewrites the two data rows of the plan-visit import sheet: row 2 gets new
# values (status/cluster/region columns were nullable and are now filled
# in), and a second data row (row 3) is added for a second outlet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 2 & 3 --------------------------------------------------------
# Write in this particular order (interleaving row 2 / row 3) so that
# newly-introduced shared strings are interned in the same sequence the
# source workbook uses.
$ws.Range("A2").Value = "CV.TOP"
$ws.Range("B2").Value = "REALME"
$ws.Range("A3").Value = "CV.TOP"
$ws.Range("B3").Value = "REALME"

$ws.Range("C2").Value = "MAINTAIN"
$ws.Range("C3").Value = "BIGCIREBON"
$ws.Range("D2").Value = "CRBASC2"
$ws.Range("D3").Value = "KRWASC3"
$ws.Range("E3").Value = "BT03852"
$ws.Range("F3").Value = "BIOHAZARD PHONE"
$ws.Range("G3").Value = "jl pasar loji desa cinta laksana tegal waru karawang"
$ws.Range("H3").Value = "KARAWANG 1"
$ws.Range("E2").Value = "BT01153"
$ws.Range("F2").Value = "GALAXY CELL"
$ws.Range("G2").Value = "Ds Rambatan Wetan Blok Pecuk Katapang Rt 29/08 (Samping SMP Hidayatul Mujahidin) Indramayu"
$ws.Range("H2").Value = "INDRAMAYU"
$ws.Range("I3").Value = "UNMAINTAIN"
$ws.Range("I2").Value = "BIGKARAWANG"

# --- Formatting ---------------------------------------------------------
# The edited cells (region/cluster/status columns that used to be blank-able)
# pick up a distinct (but still plain black) cell style in the source file,
# so give them an explicit font to split them into their own style record.
$ws.Range("C2").Font.ThemeColor = 1
$ws.Range("D2").Font.ThemeColor = 1
$ws.Range("I2").Font.ThemeColor = 1
$ws.Range("C3").Font.ThemeColor = 1
$ws.Range("I3").Font.ThemeColor = 1

# --- Selection ------------------------------------------------------------
[void]$ws.Range("C3").Select()
